# Apply the DESaster input-data-template changes:
#  - renters sheet gains a "Listed" boolean column, and 4 new renter rows
#    (copied across from the forrent_stock sheet, now also carrying a
#    landlord name and a Listed flag)
#  - owners sheet gains a "Listed" boolean column, and 4 new owner rows
#    (copied across from the forsale_stock sheet, now flagged Listed=TRUE)
#  - forrent_stock sheet gains a (currently empty) "For Sale" column
#  - selections / active sheet tweaked to match the author's last editing
#    position

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # renters
$ws2 = $wb.Worksheets.Item(2)   # owners
$ws3 = $wb.Worksheets.Item(3)   # forsale_stock
$ws4 = $wb.Worksheets.Item(4)   # forrent_stock

# ---------------------------------------------------------------------------
# New shared strings must be introduced in this exact order so they land at
# the same shared-string-table indices the workbook ends up with:
#   "For Sale" (forrent_stock header), "Listed" (renters header, reused by
#   owners header), then the four new landlord names.
# ---------------------------------------------------------------------------

# forrent_stock: new (blank-for-now) "For Sale" column header
$ws4.Cells.Item(1, 13).Value = "For Sale"

# renters: new "Listed" column header
$ws1.Cells.Item(1, 14).Value = "Listed"

# renters: existing 4 rows are not currently listed
$ws1.Cells.Item(2, 14).Value = $false
$ws1.Cells.Item(3, 14).Value = $false
$ws1.Cells.Item(4, 14).Value = $false
$ws1.Cells.Item(5, 14).Value = $false

# renters: four new rows, copied over from forrent_stock, each now with a
# landlord (column K) and marked Listed = TRUE (column N)
$ws1.Cells.Item(6, 1).Value = "Butch"
$ws1.Cells.Item(6, 2).Value = "100 Old Ave"
$ws1.Cells.Item(6, 3).Value = "Mobile Home"
$ws1.Cells.Item(6, 4).Value = 100
$ws1.Cells.Item(6, 5).Value = 1
$ws1.Cells.Item(6, 6).Value = 1
$ws1.Cells.Item(6, 7).Value = 700
$ws1.Cells.Item(6, 8).Value = 1920
$ws1.Cells.Item(6, 9).Value = 99999
$ws1.Cells.Item(6, 10).Value = "Complete"
$ws1.Cells.Item(6, 11).Value = "Greg"
$ws1.Cells.Item(6, 12).Value = 100000000
$ws1.Cells.Item(6, 13).Value = 1
$ws1.Cells.Item(6, 14).Value = $true

$ws1.Cells.Item(7, 1).Value = "Harvey"
$ws1.Cells.Item(7, 2).Value = "101 Old Ave"
$ws1.Cells.Item(7, 3).Value = "Single Family Dwelling"
$ws1.Cells.Item(7, 4).Value = 100000
$ws1.Cells.Item(7, 5).Value = 2
$ws1.Cells.Item(7, 6).Value = 2
$ws1.Cells.Item(7, 7).Value = 5000
$ws1.Cells.Item(7, 8).Value = 1920
$ws1.Cells.Item(7, 9).Value = 9999
$ws1.Cells.Item(7, 10).Value = "Complete"
$ws1.Cells.Item(7, 11).Value = "Allison"
$ws1.Cells.Item(7, 12).Value = 100000000
$ws1.Cells.Item(7, 13).Value = 1
$ws1.Cells.Item(7, 14).Value = $true

$ws1.Cells.Item(8, 1).Value = "Lee"
$ws1.Cells.Item(8, 2).Value = "102 Old Ave"
$ws1.Cells.Item(8, 3).Value = "Mobile Home"
$ws1.Cells.Item(8, 4).Value = 10
$ws1.Cells.Item(8, 5).Value = 0
$ws1.Cells.Item(8, 6).Value = 1
$ws1.Cells.Item(8, 7).Value = 250
$ws1.Cells.Item(8, 8).Value = 1960
$ws1.Cells.Item(8, 9).Value = 9999
$ws1.Cells.Item(8, 10).Value = "Complete"
$ws1.Cells.Item(8, 11).Value = "Rachel"
$ws1.Cells.Item(8, 12).Value = 100000000
$ws1.Cells.Item(8, 13).Value = 1
$ws1.Cells.Item(8, 14).Value = $true

$ws1.Cells.Item(9, 1).Value = "Carmine"
$ws1.Cells.Item(9, 2).Value = "103 Old Ave"
$ws1.Cells.Item(9, 3).Value = "Single Family Dwelling"
$ws1.Cells.Item(9, 4).Value = 2000
$ws1.Cells.Item(9, 5).Value = 3
$ws1.Cells.Item(9, 6).Value = 2
$ws1.Cells.Item(9, 7).Value = 2000
$ws1.Cells.Item(9, 8).Value = 2010
$ws1.Cells.Item(9, 9).Value = 800000
$ws1.Cells.Item(9, 10).Value = "Complete"
$ws1.Cells.Item(9, 11).Value = "Larry"
$ws1.Cells.Item(9, 12).Value = 100000000
$ws1.Cells.Item(9, 13).Value = 1
$ws1.Cells.Item(9, 14).Value = $true

# owners: new "Listed" column header (reuses the shared string created above)
$ws2.Cells.Item(1, 13).Value = "Listed"

# owners: existing 4 rows are not currently listed
$ws2.Cells.Item(2, 13).Value = $false
$ws2.Cells.Item(3, 13).Value = $false
$ws2.Cells.Item(4, 13).Value = $false
$ws2.Cells.Item(5, 13).Value = $false

# owners: four new rows, copied over from forsale_stock, each now marked
# Listed = TRUE (column M)
$ws2.Cells.Item(6, 1).Value = "Jerome"
$ws2.Cells.Item(6, 2).Value = 100000000
$ws2.Cells.Item(6, 3).Value = 1
$ws2.Cells.Item(6, 4).Value = "100 New Ave"
$ws2.Cells.Item(6, 5).Value = 1000
$ws2.Cells.Item(6, 6).Value = "Mobile Home"
$ws2.Cells.Item(6, 7).Value = 1
$ws2.Cells.Item(6, 8).Value = 1
$ws2.Cells.Item(6, 9).Value = 1100
$ws2.Cells.Item(6, 10).Value = 1920
$ws2.Cells.Item(6, 11).Value = 100000
$ws2.Cells.Item(6, 12).Value = "Slight"
$ws2.Cells.Item(6, 13).Value = $true

$ws2.Cells.Item(7, 1).Value = "Barbara"
$ws2.Cells.Item(7, 2).Value = 100000000
$ws2.Cells.Item(7, 3).Value = 1
$ws2.Cells.Item(7, 4).Value = "101 New Ave"
$ws2.Cells.Item(7, 5).Value = 4000
$ws2.Cells.Item(7, 6).Value = "Single Family Dwelling"
$ws2.Cells.Item(7, 7).Value = 4
$ws2.Cells.Item(7, 8).Value = 5
$ws2.Cells.Item(7, 9).Value = 5000
$ws2.Cells.Item(7, 10).Value = 1920
$ws2.Cells.Item(7, 11).Value = 10000000
$ws2.Cells.Item(7, 12).Value = "Extensive"
$ws2.Cells.Item(7, 13).Value = $true

$ws2.Cells.Item(8, 1).Value = "Lucius"
$ws2.Cells.Item(8, 2).Value = 100000000
$ws2.Cells.Item(8, 3).Value = 1
$ws2.Cells.Item(8, 4).Value = "102 New Ave"
$ws2.Cells.Item(8, 5).Value = 1000
$ws2.Cells.Item(8, 6).Value = "Single Family Dwelling"
$ws2.Cells.Item(8, 7).Value = 2
$ws2.Cells.Item(8, 8).Value = 1
$ws2.Cells.Item(8, 9).Value = 1200
$ws2.Cells.Item(8, 10).Value = 1960
$ws2.Cells.Item(8, 11).Value = 10000
$ws2.Cells.Item(8, 12).Value = "Moderate"
$ws2.Cells.Item(8, 13).Value = $true

$ws2.Cells.Item(9, 1).Value = "Dick"
$ws2.Cells.Item(9, 2).Value = 100000000
$ws2.Cells.Item(9, 3).Value = 1
$ws2.Cells.Item(9, 4).Value = "103 New Ave"
$ws2.Cells.Item(9, 5).Value = 2000
$ws2.Cells.Item(9, 6).Value = "Single Family Dwelling"
$ws2.Cells.Item(9, 7).Value = 3
$ws2.Cells.Item(9, 8).Value = 2
$ws2.Cells.Item(9, 9).Value = 2000
$ws2.Cells.Item(9, 10).Value = 2010
$ws2.Cells.Item(9, 11).Value = 700000
$ws2.Cells.Item(9, 12).Value = "Complete"
$ws2.Cells.Item(9, 13).Value = $true

# ---------------------------------------------------------------------------
# View state: restore each sheet's selection. The sheet that should end up
# active (owners) must be activated, and have its selection applied, last --
# otherwise selecting a range on another sheet steals the "active sheet"
# status back.
# ---------------------------------------------------------------------------

$ws1.Range("N3").Select()
$ws3.Range("A1:L5").Select()
$ws4.Range("A2:L5").Select()

$ws2.Activate()
$ws2.Range("L12").Select()
